$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A41:A53").NumberFormat = "@"
$ws.Range("A41").Value = "2026-01-28"
$ws.Range("B41").Value = "15:31:42"
$ws.Range("C41").Value = "15:00"
$ws.Range("D41").Value = "Bathroom"
$ws.Range("E41").Value = "No Motion"
$ws.Range("F41").Value = "Inactive"
$ws.Range("A42").Value = "2026-01-28"
$ws.Range("B42").Value = "15:31:44"
$ws.Range("C42").Value = "15:00"
$ws.Range("D42").Value = "Bathroom"
$ws.Range("E42").Value = "No Motion"
$ws.Range("F42").Value = "Inactive"
$ws.Range("A43").Value = "2026-01-28"
$ws.Range("B43").Value = "15:31:49"
$ws.Range("C43").Value = "15:00"
$ws.Range("D43").Value = "Bathroom"
$ws.Range("E43").Value = "No Motion"
$ws.Range("F43").Value = "Inactive"
$ws.Range("A44").Value = "2026-01-28"
$ws.Range("B44").Value = "15:31:54"
$ws.Range("C44").Value = "15:00"
$ws.Range("D44").Value = "Bathroom"
$ws.Range("E44").Value = "No Motion"
$ws.Range("F44").Value = "Inactive"
$ws.Range("A45").Value = "2026-01-28"
$ws.Range("B45").Value = "15:31:59"
$ws.Range("C45").Value = "15:00"
$ws.Range("D45").Value = "Bathroom"
$ws.Range("E45").Value = "No Motion"
$ws.Range("F45").Value = "Inactive"
$ws.Range("A46").Value = "2026-01-28"
$ws.Range("B46").Value = "15:32:04"
$ws.Range("C46").Value = "15:00"
$ws.Range("D46").Value = "Bathroom"
$ws.Range("E46").Value = "No Motion"
$ws.Range("F46").Value = "Inactive"
$ws.Range("A47").Value = "2026-01-28"
$ws.Range("B47").Value = "15:32:09"
$ws.Range("C47").Value = "15:00"
$ws.Range("D47").Value = "Bathroom"
$ws.Range("E47").Value = "No Motion"
$ws.Range("F47").Value = "Inactive"
$ws.Range("A48").Value = "2026-01-28"
$ws.Range("B48").Value = "15:32:14"
$ws.Range("C48").Value = "15:00"
$ws.Range("D48").Value = "Bathroom"
$ws.Range("E48").Value = "No Motion"
$ws.Range("F48").Value = "Inactive"
$ws.Range("A49").Value = "2026-01-28"
$ws.Range("B49").Value = "15:32:19"
$ws.Range("C49").Value = "15:00"
$ws.Range("D49").Value = "Bathroom"
$ws.Range("E49").Value = "No Motion"
$ws.Range("F49").Value = "Inactive"
$ws.Range("A50").Value = "2026-01-28"
$ws.Range("B50").Value = "15:32:24"
$ws.Range("C50").Value = "15:00"
$ws.Range("D50").Value = "Bathroom"
$ws.Range("E50").Value = "No Motion"
$ws.Range("F50").Value = "Inactive"
$ws.Range("A51").Value = "2026-01-28"
$ws.Range("B51").Value = "15:32:29"
$ws.Range("C51").Value = "15:00"
$ws.Range("D51").Value = "Bathroom"
$ws.Range("E51").Value = "No Motion"
$ws.Range("F51").Value = "Inactive"
$ws.Range("A52").Value = "2026-01-28"
$ws.Range("B52").Value = "15:32:34"
$ws.Range("C52").Value = "15:00"
$ws.Range("D52").Value = "Bathroom"
$ws.Range("E52").Value = "No Motion"
$ws.Range("F52").Value = "Inactive"
$ws.Range("A53").Value = "2026-01-28"
$ws.Range("B53").Value = "15:32:39"
$ws.Range("C53").Value = "15:00"
$ws.Range("D53").Value = "Bathroom"
$ws.Range("E53").Value = "No Motion"
$ws.Range("F53").Value = "Inactive"
$ws.Range("A41:A53").ClearFormats()

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A41:A57").NumberFormat = "@"
$ws.Range("E41:E57").NumberFormat = "@"
$ws.Range("A41").Value = "2026-01-28"
$ws.Range("B41").Value = "15:31:42"
$ws.Range("C41").Value = "15:00"
$ws.Range("D41").Value = "Bathroom"
$ws.Range("E41").Value = "88.2%"
$ws.Range("F41").Value = "Active"
$ws.Range("A42").Value = "2026-01-28"
$ws.Range("B42").Value = "15:31:42"
$ws.Range("C42").Value = "15:00"
$ws.Range("D42").Value = "Bathroom"
$ws.Range("E42").Value = "87.2%"
$ws.Range("F42").Value = "Active"
$ws.Range("A43").Value = "2026-01-28"
$ws.Range("B43").Value = "15:31:45"
$ws.Range("C43").Value = "15:00"
$ws.Range("D43").Value = "Bathroom"
$ws.Range("E43").Value = "88.2%"
$ws.Range("F43").Value = "Active"
$ws.Range("A44").Value = "2026-01-28"
$ws.Range("B44").Value = "15:31:49"
$ws.Range("C44").Value = "15:00"
$ws.Range("D44").Value = "Bathroom"
$ws.Range("E44").Value = "87.2%"
$ws.Range("F44").Value = "Active"
$ws.Range("A45").Value = "2026-01-28"
$ws.Range("B45").Value = "15:31:53"
$ws.Range("C45").Value = "15:00"
$ws.Range("D45").Value = "Bathroom"
$ws.Range("E45").Value = "88.2%"
$ws.Range("F45").Value = "Active"
$ws.Range("A46").Value = "2026-01-28"
$ws.Range("B46").Value = "15:31:57"
$ws.Range("C46").Value = "15:00"
$ws.Range("D46").Value = "Bathroom"
$ws.Range("E46").Value = "88.2%"
$ws.Range("F46").Value = "Active"
$ws.Range("A47").Value = "2026-01-28"
$ws.Range("B47").Value = "15:32:01"
$ws.Range("C47").Value = "15:00"
$ws.Range("D47").Value = "Bathroom"
$ws.Range("E47").Value = "87.3%"
$ws.Range("F47").Value = "Active"
$ws.Range("A48").Value = "2026-01-28"
$ws.Range("B48").Value = "15:32:05"
$ws.Range("C48").Value = "15:00"
$ws.Range("D48").Value = "Bathroom"
$ws.Range("E48").Value = "88.2%"
$ws.Range("F48").Value = "Active"
$ws.Range("A49").Value = "2026-01-28"
$ws.Range("B49").Value = "15:32:09"
$ws.Range("C49").Value = "15:00"
$ws.Range("D49").Value = "Bathroom"
$ws.Range("E49").Value = "87.3%"
$ws.Range("F49").Value = "Active"
$ws.Range("A50").Value = "2026-01-28"
$ws.Range("B50").Value = "15:32:13"
$ws.Range("C50").Value = "15:00"
$ws.Range("D50").Value = "Bathroom"
$ws.Range("E50").Value = "88.2%"
$ws.Range("F50").Value = "Active"
$ws.Range("A51").Value = "2026-01-28"
$ws.Range("B51").Value = "15:32:17"
$ws.Range("C51").Value = "15:00"
$ws.Range("D51").Value = "Bathroom"
$ws.Range("E51").Value = "88.2%"
$ws.Range("F51").Value = "Active"
$ws.Range("A52").Value = "2026-01-28"
$ws.Range("B52").Value = "15:32:21"
$ws.Range("C52").Value = "15:00"
$ws.Range("D52").Value = "Bathroom"
$ws.Range("E52").Value = "87.3%"
$ws.Range("F52").Value = "Active"
$ws.Range("A53").Value = "2026-01-28"
$ws.Range("B53").Value = "15:32:25"
$ws.Range("C53").Value = "15:00"
$ws.Range("D53").Value = "Bathroom"
$ws.Range("E53").Value = "88.2%"
$ws.Range("F53").Value = "Active"
$ws.Range("A54").Value = "2026-01-28"
$ws.Range("B54").Value = "15:32:30"
$ws.Range("C54").Value = "15:00"
$ws.Range("D54").Value = "Bathroom"
$ws.Range("E54").Value = "87.3%"
$ws.Range("F54").Value = "Active"
$ws.Range("A55").Value = "2026-01-28"
$ws.Range("B55").Value = "15:32:34"
$ws.Range("C55").Value = "15:00"
$ws.Range("D55").Value = "Bathroom"
$ws.Range("E55").Value = "88.2%"
$ws.Range("F55").Value = "Active"
$ws.Range("A56").Value = "2026-01-28"
$ws.Range("B56").Value = "15:32:38"
$ws.Range("C56").Value = "15:00"
$ws.Range("D56").Value = "Bathroom"
$ws.Range("E56").Value = "88.2%"
$ws.Range("F56").Value = "Active"
$ws.Range("A57").Value = "2026-01-28"
$ws.Range("B57").Value = "15:32:42"
$ws.Range("C57").Value = "15:00"
$ws.Range("D57").Value = "Bathroom"
$ws.Range("E57").Value = "87.3%"
$ws.Range("F57").Value = "Active"
$ws.Range("A41:A57").ClearFormats()
$ws.Range("E41:E57").ClearFormats()

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A41:A57").NumberFormat = "@"
$ws.Range("A41").Value = "2026-01-28"
$ws.Range("B41").Value = "15:31:42"
$ws.Range("C41").Value = "15:00"
$ws.Range("D41").Value = "Bathroom"
$ws.Range("E41").Value = "22.9C"
$ws.Range("F41").Value = "Active"
$ws.Range("A42").Value = "2026-01-28"
$ws.Range("B42").Value = "15:31:43"
$ws.Range("C42").Value = "15:00"
$ws.Range("D42").Value = "Bathroom"
$ws.Range("E42").Value = "22.9C"
$ws.Range("F42").Value = "Active"
$ws.Range("A43").Value = "2026-01-28"
$ws.Range("B43").Value = "15:31:45"
$ws.Range("C43").Value = "15:00"
$ws.Range("D43").Value = "Bathroom"
$ws.Range("E43").Value = "22.9C"
$ws.Range("F43").Value = "Active"
$ws.Range("A44").Value = "2026-01-28"
$ws.Range("B44").Value = "15:31:49"
$ws.Range("C44").Value = "15:00"
$ws.Range("D44").Value = "Bathroom"
$ws.Range("E44").Value = "22.9C"
$ws.Range("F44").Value = "Active"
$ws.Range("A45").Value = "2026-01-28"
$ws.Range("B45").Value = "15:31:53"
$ws.Range("C45").Value = "15:00"
$ws.Range("D45").Value = "Bathroom"
$ws.Range("E45").Value = "22.9C"
$ws.Range("F45").Value = "Active"
$ws.Range("A46").Value = "2026-01-28"
$ws.Range("B46").Value = "15:31:58"
$ws.Range("C46").Value = "15:00"
$ws.Range("D46").Value = "Bathroom"
$ws.Range("E46").Value = "22.9C"
$ws.Range("F46").Value = "Active"
$ws.Range("A47").Value = "2026-01-28"
$ws.Range("B47").Value = "15:32:02"
$ws.Range("C47").Value = "15:00"
$ws.Range("D47").Value = "Bathroom"
$ws.Range("E47").Value = "22.9C"
$ws.Range("F47").Value = "Active"
$ws.Range("A48").Value = "2026-01-28"
$ws.Range("B48").Value = "15:32:06"
$ws.Range("C48").Value = "15:00"
$ws.Range("D48").Value = "Bathroom"
$ws.Range("E48").Value = "22.9C"
$ws.Range("F48").Value = "Active"
$ws.Range("A49").Value = "2026-01-28"
$ws.Range("B49").Value = "15:32:10"
$ws.Range("C49").Value = "15:00"
$ws.Range("D49").Value = "Bathroom"
$ws.Range("E49").Value = "22.9C"
$ws.Range("F49").Value = "Active"
$ws.Range("A50").Value = "2026-01-28"
$ws.Range("B50").Value = "15:32:14"
$ws.Range("C50").Value = "15:00"
$ws.Range("D50").Value = "Bathroom"
$ws.Range("E50").Value = "22.9C"
$ws.Range("F50").Value = "Active"
$ws.Range("A51").Value = "2026-01-28"
$ws.Range("B51").Value = "15:32:18"
$ws.Range("C51").Value = "15:00"
$ws.Range("D51").Value = "Bathroom"
$ws.Range("E51").Value = "22.9C"
$ws.Range("F51").Value = "Active"
$ws.Range("A52").Value = "2026-01-28"
$ws.Range("B52").Value = "15:32:22"
$ws.Range("C52").Value = "15:00"
$ws.Range("D52").Value = "Bathroom"
$ws.Range("E52").Value = "22.9C"
$ws.Range("F52").Value = "Active"
$ws.Range("A53").Value = "2026-01-28"
$ws.Range("B53").Value = "15:32:26"
$ws.Range("C53").Value = "15:00"
$ws.Range("D53").Value = "Bathroom"
$ws.Range("E53").Value = "22.9C"
$ws.Range("F53").Value = "Active"
$ws.Range("A54").Value = "2026-01-28"
$ws.Range("B54").Value = "15:32:30"
$ws.Range("C54").Value = "15:00"
$ws.Range("D54").Value = "Bathroom"
$ws.Range("E54").Value = "22.9C"
$ws.Range("F54").Value = "Active"
$ws.Range("A55").Value = "2026-01-28"
$ws.Range("B55").Value = "15:32:34"
$ws.Range("C55").Value = "15:00"
$ws.Range("D55").Value = "Bathroom"
$ws.Range("E55").Value = "22.9C"
$ws.Range("F55").Value = "Active"
$ws.Range("A56").Value = "2026-01-28"
$ws.Range("B56").Value = "15:32:38"
$ws.Range("C56").Value = "15:00"
$ws.Range("D56").Value = "Bathroom"
$ws.Range("E56").Value = "22.9C"
$ws.Range("F56").Value = "Active"
$ws.Range("A57").Value = "2026-01-28"
$ws.Range("B57").Value = "15:32:42"
$ws.Range("C57").Value = "15:00"
$ws.Range("D57").Value = "Bathroom"
$ws.Range("E57").Value = "22.9C"
$ws.Range("F57").Value = "Active"
$ws.Range("A41:A57").ClearFormats()
